$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.734.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.611.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.80%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5138"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.008"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2539"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06134"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07524"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.625.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.294"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.852.90"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5352"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7772"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.809.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.551"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "181.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.880"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.24%  "
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.005"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.976"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1195"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.226"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.361"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05834"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.234"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.00%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.293"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.69%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.279"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.569"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9499"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.390"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.714"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5660"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01582"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.003"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8283"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.572"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.008.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "98.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.777.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈106"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.36%  "
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.869"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05163"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4221"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.85%  "
